$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Note: Shape.Left/Top/Width/Height round-trip through a single-precision
# (float32) COM property, then the EMU actually written is
# trunc(float32(points) * 12700). The literals below were chosen so that,
# after that float32 truncation, they land exactly on the target EMU value
# from the target OOXML (rather than using the naive EMU/12700.0 which can
# be off by a rounding ULP).

# Picture 1 - shift left (x: -968936 -> -1166155 EMU)
$shp = $s.Shapes.Item(1)
$shp.Left = -91.82323455810547

# Picture 2 - reposition & resize
#   x: 8597516 -> 8237743 EMU, y: 1163512 -> 896667 EMU
#   cx: 6354163 -> 7394279 EMU, cy: 3859167 -> 4490876 EMU
$shp = $s.Shapes.Item(2)
$shp.Left = 648.6412353515625
$shp.Top = 70.60370635986328
$shp.Width = 582.2267456054688
$shp.Height = 353.6123046875

# Picture 3 - shift left (x: 715278 -> 518059 EMU)
$shp = $s.Shapes.Item(3)
$shp.Left = 40.792049407958984

# Picture 8 - shift left (x: 8540367 -> 8343148 EMU)
$shp = $s.Shapes.Item(4)
$shp.Left = 656.9407958984375

# TextBox 18 ("a") - shift left (x: 748700 -> 551481 EMU)
$shp = $s.Shapes.Item(5)
$shp.Left = 43.423702239990234

# TextBox 19 ("b") - shift left (x: 8245088 -> 8047869 EMU)
$shp = $s.Shapes.Item(6)
$shp.Left = 633.6904907226562

# TextBox 20 ("d") - shift left (x: 8245088 -> 8047869 EMU)
$shp = $s.Shapes.Item(7)
$shp.Left = 633.6904907226562

# TextBox 21 ("c") - shift left (x: 748700 -> 551481 EMU)
$shp = $s.Shapes.Item(8)
$shp.Left = 43.423702239990234

# Slide 2: table cell (row 8, col 18) text re-set so PowerPoint marks the run as "dirty" (proofed)
$s2 = $p.Slides.Item(2)
$tbl = $s2.Shapes.Item(1).Table
$cell = $tbl.Cell(8, 18)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Text = "0"
